# Updated cryptos list on Mon Aug 28 19:07:57 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values that are really text (locale-formatted
# numbers such as "26.235.06" or "1.003" that are NOT valid Excel numbers in
# general). Force the whole data range to Text format before writing so the
# COM layer does not auto-convert digit-and-dot strings into real numbers,
# then drop back to the default "Normal" style so no stray formatting is
# left behind on the cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Map of row -> (new D value, new E value). $null means "leave unchanged".
$sub5 = [char]0x2085
$updates = @{
    2  = @("26.235.06", "  -0.36%  ")
    3  = @("1.657.53",  "  -0.62%  ")
    4  = @("1.003",     "  -0.49%  ")
    5  = @("219.46",    "  -0.43%  ")
    6  = @("0.5261",    "  -0.95%  ")
    7  = @($null,       "  -0.47%  ")
    8  = @("0.2672",    "  +0.89%  ")
    9  = @("0.06367",   "  -0.09%  ")
    10 = @("20.74",     "  -0.69%  ")
    11 = @("0.07737",   "  -1.23%  ")
    12 = @("4.605",     "  +1.77%  ")
    13 = @("1.673.08",  "  +0.14%  ")
    14 = @("1.883.66",  $null)
    15 = @("0.5644",    "  +0.79%  ")
    16 = @(("0.0{0}8243" -f $sub5), "  +0.95%  ")
    17 = @("65.50",     "  -0.51%  ")
    18 = @("26.226.82", "  -0.35%  ")
    19 = @($null,       "  -0.51%  ")
    20 = @("4.714",     "  -0.28%  ")
    21 = @("10.41",     "  +1.36%  ")
    22 = @("192.36",    "  -2.60%  ")
    23 = @("6.018",     "  -0.53%  ")
    24 = @($null,       "  -0.48%  ")
    25 = @("143.97",    "  -1.56%  ")
    26 = @("0.1208",    "  -1.09%  ")
    27 = @("7.283",     "  +0.42%  ")
    28 = @("15.96",     "  -1.41%  ")
    30 = @("0.05630",   "  -4.60%  ")
    31 = @($null,       "  -0.52%  ")
    32 = @("3.507",     "  -1.16%  ")
    33 = @("3.378",     "  +1.36%  ")
    34 = @("1.586",     "  -1.19%  ")
    35 = @("0.9539",    "  -0.77%  ")
    36 = @("2.800",     "  -1.00%  ")
    37 = @("2.409",     "  -0.93%  ")
    38 = @("0.5775",    "  -0.84%  ")
    39 = @("6.026",     "  +1.07%  ")
    40 = @($null,       "  -0.94%  ")
    41 = @("1.003",     "  -0.55%  ")
    42 = @("0.8424",    "  -1.79%  ")
    43 = @("101.96",    "  -0.94%  ")
    44 = @("1.013.58",  "  -6.02%  ")
    45 = @("1.794.87",  "  -0.62%  ")
    46 = @("58.57",     $null)
    47 = @("1.006",     "  -0.96%  ")
    48 = @("0.05335",   "  +3.56%  ")
    51 = @("0.09789",   "  +1.93%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]
    if ($null -ne $dVal) {
        $ws.Cells.Item($row, 4).Value = $dVal
    }
    if ($null -ne $eVal) {
        $ws.Cells.Item($row, 5).Value = $eVal
    }
}

# Rows 49 and 50 swap coin identity (EnergySwap <-> Mantle) and get new
# price/volume figures in the process.
$ws.Cells.Item(49, 2).Value = "Mantle"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(49, 4).Value = "0.4348"
$ws.Cells.Item(49, 5).Value = "  -1.30%  "

$ws.Cells.Item(50, 2).Value = "EnergySwap"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(50, 4).Value = "8.015"
$ws.Cells.Item(50, 5).Value = "  -0.71%  "

# Drop the transient Text number-format back to the sheet's default style so
# the cells end up with no explicit style override, matching the original
# workbook's formatting.
$priceRange.Style = "Normal"
